$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @(2, 45993, 13, 0, "02.12.202513"),
    @(3, 45993, 14, 0, "02.12.202514"),
    @(4, 45993, 15, 0.331, "02.12.202515"),
    @(5, 45993, 16, 0.234, "02.12.202516"),
    @(6, 45993, 17, 0.089, "02.12.202517"),
    @(7, 45993, 18, 0, "02.12.202518"),
    @(8, 45993, 19, 0, "02.12.202519"),
    @(9, 45993, 20, 0, "02.12.202520"),
    @(10, 45993, 21, 0, "02.12.202521"),
    @(11, 45993, 22, 0, "02.12.202522"),
    @(12, 45993, 23, 0, "02.12.202523"),
    @(13, 45993, 24, 0, "02.12.202524"),
    @(14, 45994, 1, 0, "03.12.20251"),
    @(15, 45994, 2, 0, "03.12.20252"),
    @(16, 45994, 3, 0, "03.12.20253"),
    @(17, 45994, 4, 0, "03.12.20254"),
    @(18, 45994, 5, 0, "03.12.20255"),
    @(19, 45994, 6, 0, "03.12.20256"),
    @(20, 45994, 7, 0, "03.12.20257"),
    @(21, 45994, 8, 0, "03.12.20258"),
    @(22, 45994, 9, 0, "03.12.20259"),
    @(23, 45994, 10, 0.06, "03.12.202510"),
    @(24, 45994, 11, 0.247, "03.12.202511"),
    @(25, 45994, 12, 0.613, "03.12.202512"),
    @(26, 45994, 13, 0.799, "03.12.202513"),
    @(27, 45994, 14, 0.846, "03.12.202514"),
    @(28, 45994, 15, 0.786, "03.12.202515"),
    @(29, 45994, 16, 0.344, "03.12.202516"),
    @(30, 45994, 17, 0.147, "03.12.202517"),
    @(31, 45994, 18, 0.01, "03.12.202518"),
    @(32, 45994, 19, 0, "03.12.202519"),
    @(33, 45994, 20, 0, "03.12.202520"),
    @(34, 45994, 21, 0, "03.12.202521"),
    @(35, 45994, 22, 0, "03.12.202522"),
    @(36, 45994, 23, 0, "03.12.202523"),
    @(37, 45994, 24, 0, "03.12.202524"),
    @(38, 45995, 1, 0, "04.12.20251"),
    @(39, 45995, 2, 0, "04.12.20252"),
    @(40, 45995, 3, 0, "04.12.20253"),
    @(41, 45995, 4, 0, "04.12.20254"),
    @(42, 45995, 5, 0, "04.12.20255"),
    @(43, 45995, 6, 0, "04.12.20256"),
    @(44, 45995, 7, 0, "04.12.20257"),
    @(45, 45995, 8, 0, "04.12.20258"),
    @(46, 45995, 9, 0, "04.12.20259"),
    @(47, 45995, 10, 0.068, "04.12.202510"),
    @(48, 45995, 11, 0.266, "04.12.202511"),
    @(49, 45995, 12, 0.657, "04.12.202512"),
    @(50, 45995, 13, 0.981, "04.12.202513"),
    @(51, 45995, 14, 1.051, "04.12.202514"),
    @(52, 45995, 15, 0.863, "04.12.202515"),
    @(53, 45995, 16, 0.496, "04.12.202516"),
    @(54, 45995, 17, 0.181, "04.12.202517"),
    @(55, 45995, 18, 0.018, "04.12.202518"),
    @(56, 45995, 19, 0, "04.12.202519"),
    @(57, 45995, 20, 0, "04.12.202520"),
    @(58, 45995, 21, 0, "04.12.202521"),
    @(59, 45995, 22, 0, "04.12.202522"),
    @(60, 45995, 23, 0, "04.12.202523"),
    @(61, 45995, 24, 0, "04.12.202524"),
    @(62, 45996, 1, 0, "05.12.20251"),
    @(63, 45996, 2, 0, "05.12.20252"),
    @(64, 45996, 3, 0, "05.12.20253"),
    @(65, 45996, 4, 0, "05.12.20254"),
    @(66, 45996, 5, 0, "05.12.20255"),
    @(67, 45996, 6, 0, "05.12.20256"),
    @(68, 45996, 7, 0, "05.12.20257"),
    @(69, 45996, 8, 0, "05.12.20258"),
    @(70, 45996, 9, 0, "05.12.20259"),
    @(71, 45996, 10, 0.057, "05.12.202510"),
    @(72, 45996, 11, 0.294, "05.12.202511"),
    @(73, 45996, 12, 0.603, "05.12.202512"),
    @(74, 45996, 13, 0.787, "05.12.202513"),
    @(75, 45996, 14, 0.865, "05.12.202514"),
    @(76, 45996, 15, 0.849, "05.12.202515"),
    @(77, 45996, 16, 0.472, "05.12.202516"),
    @(78, 45996, 17, 0.188, "05.12.202517"),
    @(79, 45996, 18, 0.013, "05.12.202518"),
    @(80, 45996, 19, 0, "05.12.202519"),
    @(81, 45996, 20, 0, "05.12.202520"),
    @(82, 45996, 21, 0, "05.12.202521"),
    @(83, 45996, 22, 0, "05.12.202522"),
    @(84, 45996, 23, 0, "05.12.202523"),
    @(85, 45996, 24, 0, "05.12.202524"),
    @(86, 45997, 1, 0, "06.12.20251"),
    @(87, 45997, 2, 0, "06.12.20252"),
    @(88, 45997, 3, 0, "06.12.20253"),
    @(89, 45997, 4, 0, "06.12.20254"),
    @(90, 45997, 5, 0, "06.12.20255"),
    @(91, 45997, 6, 0, "06.12.20256"),
    @(92, 45997, 7, 0, "06.12.20257"),
    @(93, 45997, 8, 0, "06.12.20258"),
    @(94, 45997, 9, 0, "06.12.20259"),
    @(95, 45997, 10, 0.07199999999999999, "06.12.202510"),
    @(96, 45997, 11, 0.262, "06.12.202511"),
    @(97, 45997, 12, 0.638, "06.12.202512"),
    @(98, 45997, 13, 1.028, "06.12.202513"),
    @(99, 45997, 14, 1.192, "06.12.202514"),
    @(100, 45997, 15, 0.86, "06.12.202515"),
    @(101, 45997, 16, 0.523, "06.12.202516"),
    @(102, 45997, 17, 0.19, "06.12.202517"),
    @(103, 45997, 18, 0.011, "06.12.202518"),
    @(104, 45997, 19, 0, "06.12.202519"),
    @(105, 45997, 20, 0, "06.12.202520"),
    @(106, 45997, 21, 0, "06.12.202521"),
    @(107, 45997, 22, 0, "06.12.202522"),
    @(108, 45997, 23, 0, "06.12.202523"),
    @(109, 45997, 24, 0, "06.12.202524"),
    @(110, 45998, 1, 0, "07.12.20251"),
    @(111, 45998, 2, 0, "07.12.20252"),
    @(112, 45998, 3, 0, "07.12.20253"),
    @(113, 45998, 4, 0, "07.12.20254"),
    @(114, 45998, 5, 0, "07.12.20255"),
    @(115, 45998, 6, 0, "07.12.20256"),
    @(116, 45998, 7, 0, "07.12.20257"),
    @(117, 45998, 8, 0, "07.12.20258"),
    @(118, 45998, 9, 0, "07.12.20259"),
    @(119, 45998, 10, 0.056, "07.12.202510"),
    @(120, 45998, 11, 0.271, "07.12.202511"),
    @(121, 45998, 12, 0.571, "07.12.202512"),
    @(122, 45998, 13, 0.799, "07.12.202513"),
    @(123, 45998, 14, 0.8070000000000001, "07.12.202514"),
    @(124, 45998, 15, 0.67, "07.12.202515"),
    @(125, 45998, 16, 0.362, "07.12.202516"),
    @(126, 45998, 17, 0.158, "07.12.202517"),
    @(127, 45998, 18, 0.011, "07.12.202518"),
    @(128, 45998, 19, 0, "07.12.202519"),
    @(129, 45998, 20, 0, "07.12.202520"),
    @(130, 45998, 21, 0, "07.12.202521"),
    @(131, 45998, 22, 0, "07.12.202522"),
    @(132, 45998, 23, 0, "07.12.202523"),
    @(133, 45998, 24, 0, "07.12.202524"),
    @(134, 45999, 1, 0, "08.12.20251"),
    @(135, 45999, 2, 0, "08.12.20252"),
    @(136, 45999, 3, 0, "08.12.20253"),
    @(137, 45999, 4, 0, "08.12.20254"),
    @(138, 45999, 5, 0, "08.12.20255"),
    @(139, 45999, 6, 0, "08.12.20256"),
    @(140, 45999, 7, 0, "08.12.20257"),
    @(141, 45999, 8, 0, "08.12.20258"),
    @(142, 45999, 9, 0, "08.12.20259"),
    @(143, 45999, 10, 0.033, "08.12.202510"),
    @(144, 45999, 11, 0.248, "08.12.202511"),
    @(145, 45999, 12, 0.389, "08.12.202512"),
    @(146, 45999, 13, 0.54, "08.12.202513"),
    @(147, 45999, 14, 0.499, "08.12.202514"),
    @(148, 45999, 15, 0.37, "08.12.202515"),
    @(149, 45999, 16, 0.287, "08.12.202516"),
    @(150, 45999, 17, 0.078, "08.12.202517"),
    @(151, 45999, 18, 0, "08.12.202518"),
    @(152, 45999, 19, 0, "08.12.202519"),
    @(153, 45999, 20, 0, "08.12.202520"),
    @(154, 45999, 21, 0, "08.12.202521"),
    @(155, 45999, 22, 0, "08.12.202522"),
    @(156, 45999, 23, 0, "08.12.202523"),
    @(157, 45999, 24, 0, "08.12.202524"),
    @(158, 46000, 1, 0, "09.12.20251"),
    @(159, 46000, 2, 0, "09.12.20252"),
    @(160, 46000, 3, 0, "09.12.20253"),
    @(161, 46000, 4, 0, "09.12.20254"),
    @(162, 46000, 5, 0, "09.12.20255"),
    @(163, 46000, 6, 0, "09.12.20256"),
    @(164, 46000, 7, 0, "09.12.20257"),
    @(165, 46000, 8, 0, "09.12.20258"),
    @(166, 46000, 9, 0, "09.12.20259"),
    @(167, 46000, 10, 0.018, "09.12.202510"),
    @(168, 46000, 11, 0.102, "09.12.202511"),
    @(169, 46000, 12, 0.225, "09.12.202512"),
    @(170, 46000, 13, 0.323, "09.12.202513")
)

foreach ($row in $rowsData) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
